# Scheduled-runner refresh: recompute Leve crafting-profit figures
# (currentAveragePrice / NQ / HQ price & profit columns) on every class
# worksheet using the latest market-board snapshot.

$wb = $excel.ActiveWorkbook
$mismatchCount = 0

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 2 (Leve Item ID 5489)
if ($ws.Range("G2").Value2 -ne 5489) { $mismatchCount++; Write-Host "WARNING: G2 on ALC expected 5489 got" $ws.Range("G2").Value2 }
$ws.Range("H2").Value2 = 413.33334
$ws.Range("I2").Value2 = 413.33334
$ws.Range("K2").Value2 = 413.33334
$ws.Range("M2").Value2 = -300.33334
# row 41 (Leve Item ID 5478)
if ($ws.Range("G41").Value2 -ne 5478) { $mismatchCount++; Write-Host "WARNING: G41 on ALC expected 5478 got" $ws.Range("G41").Value2 }
$ws.Range("H41").Value2 = 1564.1538
$ws.Range("I41").Value2 = 1939.1666
$ws.Range("J41").Value2 = 1242.7142
$ws.Range("K41").Value2 = 1939.1666
$ws.Range("L41").Value2 = 1242.7142
$ws.Range("M41").Value2 = -1499.1666
$ws.Range("N41").Value2 = -2122.7142
# row 51 (Leve Item ID 5486)
if ($ws.Range("G51").Value2 -ne 5486) { $mismatchCount++; Write-Host "WARNING: G51 on ALC expected 5486 got" $ws.Range("G51").Value2 }
$ws.Range("H51").Value2 = 6037.5
$ws.Range("I51").Value2 = 6560
$ws.Range("J51").Value2 = 5166.6665
$ws.Range("K51").Value2 = 6560
$ws.Range("L51").Value2 = 5166.6665
$ws.Range("M51").Value2 = -6076
$ws.Range("N51").Value2 = -6134.6665
# row 62 (Leve Item ID 27781)
if ($ws.Range("G62").Value2 -ne 27781) { $mismatchCount++; Write-Host "WARNING: G62 on ALC expected 27781 got" $ws.Range("G62").Value2 }
$ws.Range("H62").Value2 = 1977.5555
$ws.Range("I62").Value2 = 1966.3334
$ws.Range("K62").Value2 = 1966.3334
$ws.Range("M62").Value2 = -1342.3334
# row 65 (Leve Item ID 27781)
if ($ws.Range("G65").Value2 -ne 27781) { $mismatchCount++; Write-Host "WARNING: G65 on ALC expected 27781 got" $ws.Range("G65").Value2 }
$ws.Range("H65").Value2 = 1977.5555
$ws.Range("I65").Value2 = 1966.3334
$ws.Range("K65").Value2 = 9831.666999999999
$ws.Range("M65").Value2 = -6711.666999999999
# row 88 (Leve Item ID 12608)
if ($ws.Range("G88").Value2 -ne 12608) { $mismatchCount++; Write-Host "WARNING: G88 on ALC expected 12608 got" $ws.Range("G88").Value2 }
$ws.Range("H88").Value2 = 1824.6666
$ws.Range("I88").Value2 = 1650
$ws.Range("J88").Value2 = 1999.3334
$ws.Range("K88").Value2 = 1650
$ws.Range("L88").Value2 = 1999.3334
$ws.Range("M88").Value2 = -1244
$ws.Range("N88").Value2 = -2811.3334
# row 91 (Leve Item ID 12608)
if ($ws.Range("G91").Value2 -ne 12608) { $mismatchCount++; Write-Host "WARNING: G91 on ALC expected 12608 got" $ws.Range("G91").Value2 }
$ws.Range("H91").Value2 = 1824.6666
$ws.Range("I91").Value2 = 1650
$ws.Range("J91").Value2 = 1999.3334
$ws.Range("K91").Value2 = 1650
$ws.Range("L91").Value2 = 1999.3334
$ws.Range("M91").Value2 = -246
$ws.Range("N91").Value2 = -4807.3334
# row 98 (Leve Item ID 36237)
if ($ws.Range("G98").Value2 -ne 36237) { $mismatchCount++; Write-Host "WARNING: G98 on ALC expected 36237 got" $ws.Range("G98").Value2 }
$ws.Range("H98").Value2 = 2603.1428
$ws.Range("I98").Value2 = 2147.52
$ws.Range("J98").Value2 = 6400
$ws.Range("K98").Value2 = 2147.52
$ws.Range("L98").Value2 = 6400
$ws.Range("M98").Value2 = -649.52
$ws.Range("N98").Value2 = -9396
# row 111 (Leve Item ID 27768)
if ($ws.Range("G111").Value2 -ne 27768) { $mismatchCount++; Write-Host "WARNING: G111 on ALC expected 27768 got" $ws.Range("G111").Value2 }
$ws.Range("H111").Value2 = 50001250
$ws.Range("I111").Value2 = 100000000
$ws.Range("J111").Value2 = 2499
$ws.Range("K111").Value2 = 300000000
$ws.Range("L111").Value2 = 7497
$ws.Range("M111").Value2 = -299996933
$ws.Range("N111").Value2 = -13631
# row 116 (Leve Item ID 27778)
if ($ws.Range("G116").Value2 -ne 27778) { $mismatchCount++; Write-Host "WARNING: G116 on ALC expected 27778 got" $ws.Range("G116").Value2 }
$ws.Range("H116").Value2 = 11192.214
$ws.Range("I116").Value2 = 27122.75
$ws.Range("J116").Value2 = 4820
$ws.Range("K116").Value2 = 27122.75
$ws.Range("L116").Value2 = 4820
$ws.Range("M116").Value2 = -23680.75
$ws.Range("N116").Value2 = -11704
# row 121 (Leve Item ID 39731)
if ($ws.Range("G121").Value2 -ne 39731) { $mismatchCount++; Write-Host "WARNING: G121 on ALC expected 39731 got" $ws.Range("G121").Value2 }
$ws.Range("H121").Value2 = 1136.8
$ws.Range("J121").Value2 = 1371
$ws.Range("L121").Value2 = 4113
$ws.Range("N121").Value2 = -7607
# row 122 (Leve Item ID 36237)
if ($ws.Range("G122").Value2 -ne 36237) { $mismatchCount++; Write-Host "WARNING: G122 on ALC expected 36237 got" $ws.Range("G122").Value2 }
$ws.Range("H122").Value2 = 2603.1428
$ws.Range("I122").Value2 = 2147.52
$ws.Range("J122").Value2 = 6400
$ws.Range("K122").Value2 = 6442.559999999999
$ws.Range("L122").Value2 = 19200
$ws.Range("M122").Value2 = -3992.559999999999
$ws.Range("N122").Value2 = -24100
# row 138 (Leve Item ID 44169)
if ($ws.Range("G138").Value2 -ne 44169) { $mismatchCount++; Write-Host "WARNING: G138 on ALC expected 44169 got" $ws.Range("G138").Value2 }
$ws.Range("H138").Value2 = 2187.3381
$ws.Range("I138").Value2 = 1972.762
$ws.Range("J138").Value2 = 2498.1035
$ws.Range("K138").Value2 = 5918.286
$ws.Range("L138").Value2 = 7494.310500000001
$ws.Range("M138").Value2 = -778.2860000000001
$ws.Range("N138").Value2 = -17774.3105

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32 (Leve Item ID 44147)
if ($ws.Range("G32").Value2 -ne 44147) { $mismatchCount++; Write-Host "WARNING: G32 on ARM expected 44147 got" $ws.Range("G32").Value2 }
$ws.Range("H32").Value2 = 5259.961
$ws.Range("I32").Value2 = 4516.6523
$ws.Range("J32").Value2 = 12098.4
$ws.Range("K32").Value2 = 4516.6523
$ws.Range("L32").Value2 = 12098.4
$ws.Range("M32").Value2 = -4229.6523
$ws.Range("N32").Value2 = -12672.4
# row 97 (Leve Item ID 19941)
if ($ws.Range("G97").Value2 -ne 19941) { $mismatchCount++; Write-Host "WARNING: G97 on ARM expected 19941 got" $ws.Range("G97").Value2 }
$ws.Range("H97").Value2 = 630.9655
$ws.Range("I97").Value2 = 653.0417
$ws.Range("J97").Value2 = 525
$ws.Range("K97").Value2 = 653.0417
$ws.Range("L97").Value2 = 525
$ws.Range("M97").Value2 = -157.0417
$ws.Range("N97").Value2 = -1517

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 105 (Leve Item ID 19947)
if ($ws.Range("G105").Value2 -ne 19947) { $mismatchCount++; Write-Host "WARNING: G105 on BSM expected 19947 got" $ws.Range("G105").Value2 }
$ws.Range("H105").Value2 = 2384
$ws.Range("I105").Value2 = 2361.2632
$ws.Range("J105").Value2 = 2600
$ws.Range("K105").Value2 = 2361.2632
$ws.Range("L105").Value2 = 2600
$ws.Range("M105").Value2 = -614.2631999999999
$ws.Range("N105").Value2 = -6094

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31 (Leve Item ID 44023)
if ($ws.Range("G31").Value2 -ne 44023) { $mismatchCount++; Write-Host "WARNING: G31 on CRP expected 44023 got" $ws.Range("G31").Value2 }
$ws.Range("H31").Value2 = 1500.1852
$ws.Range("I31").Value2 = 1061.7142
$ws.Range("J31").Value2 = 3034.8333
$ws.Range("K31").Value2 = 1061.7142
$ws.Range("L31").Value2 = 3034.8333
$ws.Range("M31").Value2 = -766.7141999999999
$ws.Range("N31").Value2 = -3624.8333
# row 34 (Leve Item ID 44023)
if ($ws.Range("G34").Value2 -ne 44023) { $mismatchCount++; Write-Host "WARNING: G34 on CRP expected 44023 got" $ws.Range("G34").Value2 }
$ws.Range("H34").Value2 = 1500.1852
$ws.Range("I34").Value2 = 1061.7142
$ws.Range("J34").Value2 = 3034.8333
$ws.Range("K34").Value2 = 1061.7142
$ws.Range("L34").Value2 = 3034.8333
$ws.Range("M34").Value2 = -859.7141999999999
$ws.Range("N34").Value2 = -3438.8333
# row 134 (Leve Item ID 44020)
if ($ws.Range("G134").Value2 -ne 44020) { $mismatchCount++; Write-Host "WARNING: G134 on CRP expected 44020 got" $ws.Range("G134").Value2 }
$ws.Range("H134").Value2 = 2324.4666
$ws.Range("I134").Value2 = 1992.3043
$ws.Range("K134").Value2 = 5976.9129
$ws.Range("M134").Value2 = -3441.9129

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 80 (Leve Item ID 12890)
if ($ws.Range("G80").Value2 -ne 12890) { $mismatchCount++; Write-Host "WARNING: G80 on CUL expected 12890 got" $ws.Range("G80").Value2 }
$ws.Range("H80").Value2 = 1408.75
$ws.Range("I80").Value2 = 1124.5
$ws.Range("J80").Value2 = 1693
$ws.Range("K80").Value2 = 3373.5
$ws.Range("L80").Value2 = 5079
$ws.Range("M80").Value2 = -2437.5
$ws.Range("N80").Value2 = -6951
# row 83 (Leve Item ID 12890)
if ($ws.Range("G83").Value2 -ne 12890) { $mismatchCount++; Write-Host "WARNING: G83 on CUL expected 12890 got" $ws.Range("G83").Value2 }
$ws.Range("H83").Value2 = 1408.75
$ws.Range("I83").Value2 = 1124.5
$ws.Range("J83").Value2 = 1693
$ws.Range("K83").Value2 = 10120.5
$ws.Range("L83").Value2 = 15237
$ws.Range("M83").Value2 = -5440.5
$ws.Range("N83").Value2 = -24597
# row 131 (Leve Item ID 36060)
if ($ws.Range("G131").Value2 -ne 36060) { $mismatchCount++; Write-Host "WARNING: G131 on CUL expected 36060 got" $ws.Range("G131").Value2 }
$ws.Range("H131").Value2 = 14981.984
$ws.Range("I131").Value2 = 492.7
$ws.Range("J131").Value2 = 17715.81
$ws.Range("K131").Value2 = 1478.1
$ws.Range("L131").Value2 = 53147.43000000001
$ws.Range("M131").Value2 = 3561.9
$ws.Range("N131").Value2 = -63227.43000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 97 (Leve Item ID 19940)
if ($ws.Range("G97").Value2 -ne 19940) { $mismatchCount++; Write-Host "WARNING: G97 on GSM expected 19940 got" $ws.Range("G97").Value2 }
$ws.Range("H97").Value2 = 1406.35
$ws.Range("I97").Value2 = 1345.9333
$ws.Range("J97").Value2 = 1587.6
$ws.Range("K97").Value2 = 1345.9333
$ws.Range("L97").Value2 = 1587.6
$ws.Range("M97").Value2 = -849.9332999999999
$ws.Range("N97").Value2 = -2579.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 46 (Leve Item ID 5282)
if ($ws.Range("G46").Value2 -ne 5282) { $mismatchCount++; Write-Host "WARNING: G46 on LTW expected 5282 got" $ws.Range("G46").Value2 }
$ws.Range("H46").Value2 = 1284.4546
$ws.Range("I46").Value2 = 702
$ws.Range("J46").Value2 = 2303.75
$ws.Range("K46").Value2 = 702
$ws.Range("L46").Value2 = 2303.75
$ws.Range("M46").Value2 = -514
$ws.Range("N46").Value2 = -2679.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 100 (Leve Item ID 19981)
if ($ws.Range("G100").Value2 -ne 19981) { $mismatchCount++; Write-Host "WARNING: G100 on WVR expected 19981 got" $ws.Range("G100").Value2 }
$ws.Range("H100").Value2 = 719.44446
$ws.Range("I100").Value2 = 496.42856
$ws.Range("J100").Value2 = 1500
$ws.Range("K100").Value2 = 992.85712
$ws.Range("L100").Value2 = 3000
$ws.Range("M100").Value2 = -451.85712
$ws.Range("N100").Value2 = -4082
# row 122 (Leve Item ID 36208)
if ($ws.Range("G122").Value2 -ne 36208) { $mismatchCount++; Write-Host "WARNING: G122 on WVR expected 36208 got" $ws.Range("G122").Value2 }
$ws.Range("H122").Value2 = 53341.867
$ws.Range("I122").Value2 = 59072.074
$ws.Range("J122").Value2 = 1770
$ws.Range("K122").Value2 = 177216.222
$ws.Range("L122").Value2 = 5310
$ws.Range("M122").Value2 = -174766.222
$ws.Range("N122").Value2 = -10210

Write-Host "Update complete. Mismatches:" $mismatchCount